# Auto-update gym prices
$wb = $excel.ActiveWorkbook

# --- Sheet "4x4 Squat Racks": update price for "The Corporate Rack" ---
$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")
$ws1.Range("C2").Value = "$2,139.00"

# --- Sheet "Squat Stands": update price for "The Associate Squat Stand" ---
$ws2 = $wb.Worksheets.Item("Squat Stands")
$ws2.Range("C2").Value = "$1,547.00"

# --- Sheet "Leg Extensions": remove the "Selectorized Seated Leg Curl/Extension"
#     row (row 5, Stray Dog Strength), shifting the Sorinex row (row 6) up to row 5 ---
$ws3 = $wb.Worksheets.Item("Leg Extensions")
$ws3.Range("A5:F5").EntireRow.Delete()
